# Using Prodigy for Named Entity Annotation.docx
# Commit: "Added new samples from the AACT DB filtered for non-drug,
# biological and dietary supplement."
#
# Two textual edits:
#  1) Collapse "the name of the local prodigy model used to identify the
#     project, " + "e.g." + " " (three runs, the middle one wrapped in
#     gramStart/gramEnd proofErr markers) into a single run reading
#     "...the project, e.g. ".
#  2) In the --label list, replace the leading
#     "SYSTEMATIC,TRIVIAL,VACCINE,IDENTIFIER,FORMULA,ABBREV,FAMILY,MULTIPLE"
#     with "DRUG".

$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "the name of the local prodigy model used to identify the project, e.g. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "the name of the local prodigy model used to identify the project, e.g. ",
    2)
if (-not $found1) {
    throw "Could not find the 'the name of the local prodigy model...' text"
}

# --- Change 2 -------------------------------------------------------------
$labelRange = $d.Content
$found2 = $labelRange.Find.Execute(
    "SYSTEMATIC,TRIVIAL,VACCINE,IDENTIFIER,FORMULA,ABBREV,FAMILY,MULTIPLE",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the SYSTEMATIC,... label list prefix"
}

# Replace the matched range's text with "DRUG". Toggling a boolean
# character-formatting property on/off afterwards forces Word to keep the
# freshly written text ("DRUG") in its own run instead of re-merging it
# with the (identically formatted) text that follows, matching how the
# document looks after a real type-over edit.
$labelRange.Text = "DRUG"
$labelRange.Font.Bold = 1
$labelRange.Font.Bold = 0
